# Rewrite the document body's single run to the new Gutachten outline
# ("Erlassen von Verwaltungsakten" / "Belastende Verwaltungsakte" / ...).
# The body is one paragraph holding one run made up of <w:t> pieces
# separated by manual line breaks (<w:br/>); every transition in the
# source document is a *blank* line, i.e. two breaks in a row. Word
# represents a manual line break as the vertical-tab character, which
# PowerShell writes as the escape sequence `v inside a double-quoted
# string, so "`v`v" below reproduces <w:br/><w:br/>.
$d = $word.ActiveDocument

$lines = @(
    "Erlassen von Verwaltungsakten",
    "Belastende Verwaltungsakte",
    "Fall:",
    "Das Landratsamt möchte eine Anordnung erlassen, dass das Fachwerkhaus mit Biberschwanz-Dachziegeln repariert werden muss, da es durch einen Sturm beschädigt wurde und Regen durch das Dach eindringt.",
    "Die Eigentümer, Forstrat Franz Konrad und sein Bruder Georg Konrad, wollen jedoch keine Investitionen mehr tätigen und das Dach nicht reparieren, da das Haus demnächst verkauft werden soll.",
    "Das Landratsamt möchte dennoch einschreiten, da das Fachwerkhaus ein Kulturdenkmal ist und eine Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erfüllen.",
    "Die Reparatur mit Biberschwanz-Dachziegeln würde etwa 1.200 Euro kosten.",
    "Gutachten:",
    "Rechtsgrundlage:",
    "Die Rechtsgrundlage könnte §1 Abs.",
    "1 in Verbindung mit § 7 Abs.1 Satz 1 DSchG sein.",
    "Materielle Voraussetzung:",
    "Tatbestandsvoraussetzung:",
    "Kulturdenkmal: Das Fachwerkhaus müsste ein Kulturdenkmal sein.",
    "Nach § 2 Abs.",
    "1 DSchG ist eine Sache, deren Erhaltung aus heimatgeschichtlichen Gründen ein öffentliches Interesse besteht, ein Kulturdenkmal.",
    "Gefährdung: Eine Gefährdung des Kulturdenkmals durch das beschädigte Dach könnte vorliegen.",
    "Rechtsfolgenseite:",
    "Der Pflichtige: Als Pflichtige kommen sowohl Forstrat Franz Konrad als auch sein Bruder Georg Konrad in Betracht.",
    "Forstrat Franz Konrad könnte pflichtig sein, sinngemäß § 7 Abs.",
    "1 Satz 1 DSchG und § 7 PolG, da er Eigentümer einer Sache ist, von deren Zustand eine Gefahr ausgeht.",
    "Georg Konrad ist ebenfalls Eigentümer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.",
    "Letztlich ist Forstrat Franz Konrad der richtige Pflichtige, da er über ein höheres Einkommen verfügt und somit leistungsfähiger ist.",
    "Ermessen:",
    "Die Denkmalschutzbehörde hat gemäß § 7 DSchG ein Ermessen, das nach § 40 LVwVfG ausgeübt wird.",
    "Die Anordnung der BSD ist verhältnismäßig und ermessensgerecht, da eine kostengünstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erfüllen.",
    "Der Vorteil für die Allgemeinheit durch die Ansehnlichkeit des Denkmals rechtfertigt den finanziellen Nachteil für den Eigentümer (Forstrat Franz Konrad).",
    "Unmöglichkeit:",
    "Es könnte eine privatrechtliche Unmöglichkeit vorliegen, wenn das private Recht eines Dritten der Ausführung des Verwaltungsaktes entgegensteht.",
    "Hier könnte §2038 Abs.1 Satz 1 BGB relevant sein, da Georg Konrad als Miterbe allen Maßnahmen zustimmen müsste, dies jedoch verweigert.",
    "Eine Ausnahme besteht gemäß §2038 Abs.1 Satz 2 Halbsatz 2 BGB, wenn die BSD als notwendige Erhaltungsmaßnahme anzusehen ist.",
    "Angesichts des höheren Verkaufswerts des Fachwerkhauses erscheint die Maßnahme wirtschaftlich vernünftig, und Forstrat Franz Konrad kann ohne Georg Konrad handeln, wodurch keine privatrechtliche Unmöglichkeit besteht.",
    "Bestimmtheit:",
    "Nach § 37 LVwVfG muss die Anordnung bestimmt genug formuliert werden.",
    "Formelle Voraussetzung:",
    "Zuständigkeit:",
    "Sachliche Zuständigkeit: Nach §§ 7 Abs.",
    "4, 3 Abs.",
    "3, Abs.",
    "1 Nr.",
    "3 DSchG und 46 Abs.",
    "2 LBO und § 15 LVG ist das Landratsamt sachlich zuständig.",
    "Örtliche Zuständigkeit: Örtlich zuständig ist das Landratsamt Ortenaukreis gemäß § 3 Abs.1 Nr.1 LVwVfG.",
    "Verfahren:",
    "Beteiligte: Nach §§ 13 Abs.",
    "1 Nr.",
    "2 und 4 LVwVfG sind Forstrat Franz Konrad (Nr.",
    "2) und Georg Konrad (Nr.",
    "4) beteiligte.",
    "Georg Konrad ist beteiligt, da er Eigentümer des Fachwerkhauses ist (§903 BGB).",
    "Es besteht also ein rechtliches Interesse nach § 13 Abs.",
    "2 LVwVfG.",
    "Ausgeschlossene Personen/Befangenheit: Wegen der kritischen Leserbriefe des Forstrat Franz Konrad könnte eine Problematik im Sinne des § 21 LVwVfg gegeben sein.",
    "Da aber der Mitarbeiter, der den Fall bearbeitet, nach Sachverhalt nicht bekannt ist, kann dem letztlich nicht weiter nachgegangen werden.",
    "Beteiligung anderer Behörden: Nach §3 Abs.",
    "4 DSchG muss das Landesamt für Denkmalpflege angehört werden.",
    "Anhörung: Nach § 28 Abs.",
    "1 LVwVfG ist Forstrat Franz Konrad und Georg Konrad die Gelegenheit zur Äußerung zu geben.",
    "Form:",
    "Formwahl: Nach § 37 Abs.",
    "2 LVwVfG kann die Reparaturanordnung hier schriftlich erfolgen.",
    "Begründungspflicht: Nach §39 Abs.",
    "1 LVwVfG ist die schriftliche Reparaturanordnung auch schriftlich zu begründen.",
    "Rechtsbehelfsbelehrung: Nach § 37 Abs.",
    "6 LVwVfG ist eine Rechtsbehelfsbelehrung beizufügen.",
    "Bekanntgabe: Nach §43 Abs.",
    "1 LVwVfG wird ein Verwaltungsakt durch Bekanntgabe wirksam.",
    "Dem Forstrat Franz Konrad sollte der Bescheid mittels PZU nach § 3 LVwZG zugestellt und damit bekanntgegeben werden."
)

$d.Content.Text = [string]::Join("`v`v", $lines)
